$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-23 Sunday" "2025-03-24 Monday"

Replace-Text "541÷7=" "639÷2="
Replace-Text "694÷4=" "108÷2="
Replace-Text "680÷3=" "378÷9="
Replace-Text "340÷7=" "519÷3="
Replace-Text "228÷3=" "818÷8="
Replace-Text "741÷9=" "206÷4="
Replace-Text "825÷4=" "235÷2="
Replace-Text "738÷3=" "271÷5="
Replace-Text "463÷4=" "330÷4="
Replace-Text "877÷2=" "798÷7="
Replace-Text "571÷3=" "556÷4="
Replace-Text "560÷4=" "572÷8="
Replace-Text "711÷2=" "200÷3="
Replace-Text "193÷5=" "981÷7="
Replace-Text "663÷8=" "315÷8="
Replace-Text "751÷4=" "808÷6="
Replace-Text "498÷7=" "980÷6="
Replace-Text "894÷7=" "857÷8="
Replace-Text "432÷2=" "551÷6="
Replace-Text "310÷4=" "489÷8="
Replace-Text "668÷5=" "902÷7="
Replace-Text "428÷8=" "534÷4="
Replace-Text "717÷6=" "186÷2="
Replace-Text "653÷4=" "393÷6="
Replace-Text "287÷9=" "560÷7="

Write-Host "Done applying replacements"
